# Insert a new data row at row 36 (pushing existing rows 36-77 down to 37-78),
# then populate the new row 36 with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(36).Insert()

$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 44658
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = 100112021
$ws.Cells.Item(36, 7).Value = "Ají"
$ws.Cells.Item(36, 8).Value = "Americana (o)"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 80
$ws.Cells.Item(36, 11).Value = 8500
$ws.Cells.Item(36, 12).Value = 9000
$ws.Cells.Item(36, 13).Value = 8750
$ws.Cells.Item(36, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(36, 15).Value = "Región del Maule"
$ws.Cells.Item(36, 16).Value = 583
$ws.Cells.Item(36, 17).Value = 15
$ws.Cells.Item(36, 18).Value = "Hortaliza"
